# Loan RBI, Variable Instalments
# Insert a new (blank) column before column N on the "Repayment schedule"
# sheet - this shifts the old "Late"/heading/"Outstanding" columns one
# place to the right (N->O, O->P, P->Q) and leaves the new column N blank.
# Also updates the active sheet/selection to match (Repayment schedule
# becomes the selected tab, with R7 selected).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Match the new column's width to its left neighbour (column M) before
# inserting, so the inserted column N keeps a sane custom width.
$mWidth = $ws.Columns("M:M").ColumnWidth
$ws.Columns("N:N").Insert()
$ws.Columns("N:N").ColumnWidth = $mWidth

# Make "Repayment schedule" the active/selected sheet with R7 selected.
$ws.Activate()
$ws.Range("R7").Select()
